$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Capture the existing "Y" style (s=22, currently only used by D72:D75)
# before we repaint D72:D75, so the new D76 cell can keep using it. ---
$ws.Range("D72").Copy()
$ws.Range("D76").PasteSpecial(-4122)   # xlPasteFormats

# --- Reformat D72:D75 to match D2's style (s=22 -> s=18), as in the diff ---
$ws.Range("D2").Copy()
$ws.Range("D72:D75").PasteSpecial(-4122)   # xlPasteFormats

# --- New rows 76 & 77: enter values in the same order the author did
# (both TCIDs first, then row 76's remaining cells, then row 77's) so the
# shared-string table ends up in the same append order as the diff. ---
$ws.Range("A76").Value = "Authoring75"
$ws.Range("A77").Value = "Authoring76"

$ws.Range("B76").Value = "OPQA-1195|OPQA-1313|OPQA-1312|OPQA-1090|OPQA-1201"
$ws.Range("C76").Value = "Verify saving post as draft, accessing it for edit from profile,delete post from prfile"
$ws.Range("D76").Value = "Y"
$ws.Range("E76").Value = ""

$ws.Range("B77").Value = "OPQA-1196|OPQA-1200|OPQA-1199"
$ws.Range("C77").Value = "Verify draft title,access and edit draft post from post modal, delete post from post modal"
$ws.Range("D77").Value = "Y"
$ws.Range("E77").Value = ""

# --- Match cell formatting (borders/fonts) for the new rows ---
# Row 76: A=1, B=8, C=1, D=22 (already set above via copy), E=1
$ws.Range("C67").Copy()
$ws.Range("A76").PasteSpecial(-4122)
$ws.Range("C76").PasteSpecial(-4122)
$ws.Range("E76").PasteSpecial(-4122)

$ws.Range("B64").Copy()
$ws.Range("B76").PasteSpecial(-4122)

# Row 77: A=1, B=1, C=1, D=1, E=1
$ws.Range("C67").Copy()
$ws.Range("A77").PasteSpecial(-4122)
$ws.Range("B77").PasteSpecial(-4122)
$ws.Range("C77").PasteSpecial(-4122)
$ws.Range("D77").PasteSpecial(-4122)
$ws.Range("E77").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- View state: scroll + selection, matching the diff ---
$win = $excel.ActiveWindow
$win.ScrollRow = 48
$win.ScrollColumn = 1
$ws.Range("D2:D75").Select()

Write-Host "done"
